# Applies the "Add files via upload" revision to Legenda.xlsx:
#  - Scenari sheet: fill in C4 ("-"), replace the placeholder "-" values in
#    C11/C17 with the real shock descriptions, and make Scenari the
#    selected/active sheet (selection -> C5).
#  - EGQ sheet: no longer the active sheet; leave its own selection at F21.

$wb = $excel.ActiveWorkbook

$egq = $wb.Worksheets.Item("EGQ")
$scenari = $wb.Worksheets.Item("Scenari")

# --- Scenari!C4: was blank, now a literal "-" (text, quote-prefixed just
#     like the other dash placeholders on this sheet) ---
$scenari.Range("C4").Value = "'-"

# --- Scenari!C11 / C17: replace placeholder "-" text with the real values ---
$scenari.Range("C11").Value = "'ITL 10Y ITL_10yr - Spread Yield: 75.00bps"
$scenari.Range("C17").Value = "'MSCI World Net TR MSCIWLDNET - Index Level -10.00pct"

# --- Selection on EGQ moves to F21 (sheet is no longer the active tab) ---
[void]$egq.Range("F21").Select()

# --- Scenari becomes the active sheet, selection at C5 ---
$scenari.Activate()
[void]$scenari.Range("C5").Select()
